# Applies a reshuffle of rows 3-14 and 16 (columns D, J, K, L, M, O, P)
# on the active worksheet, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each destination row, for columns D,J,K,L,M,O,P.
# (taken from the source rows described by the diff's permutation)
$rowsData = @{
    3  = @{ D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
    4  = @{ D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";   P = 551 }
    5  = @{ D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";   P = 484 }
    6  = @{ D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";   P = 622 }
    7  = @{ D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    8  = @{ D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";   P = 485 }
    9  = @{ D = 44957; J = 70;  K = 1500; L = 2000; M = 1857; O = "Región Metropolitana";   P = 310 }
    10 = @{ D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    11 = @{ D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";   P = 463 }
    12 = @{ D = 44876; J = 80;  K = 6500; L = 7000; M = 6812; O = "Región Metropolitana";   P = 1135 }
    13 = @{ D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";   P = 548 }
    14 = @{ D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";   P = 544 }
    16 = @{ D = 44987; J = 130; K = 4500; L = 5000; M = 4692; O = "Región Metropolitana";   P = 782 }
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]

    $ws.Cells.Item($r, 4).Value2 = $vals.D    # D: Fecha
    $ws.Cells.Item($r, 10).Value2 = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value2 = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value2 = $vals.O   # O: Origen
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # P: Precio $/Kg
}
